$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
  @("ECs","Spn","Siglec1","ECs",1,0.3333333333333333,0.01649366666666667,0.049481,0.005209935740510855,0.005209935740510855,2,0.6666666666666666,0.216469,0.6494070000000001,0.005147493274683087,0.005147493274683087,0.003570367529666667,0.032133307767,0.00002681810918581067,0.00002681810918581067),
  @("ECs","Spn","Siglec1","MuSCs",1,0.3333333333333333,0.01649366666666667,0.049481,0.005209935740510855,0.005209935740510855,2,0.6666666666666666,0.03281033333333333,0.09843099999999999,0.0007802085757011102,0.0007802085757011101,0.0005411627012222222,0.004870464310999999,0.000004064836543598283,0.000004064836543598282),
  @("ECs","Spn","Siglec1","Resolving-Mac",1,0.3333333333333333,0.01649366666666667,0.049481,0.005209935740510855,0.005209935740510855,3,1,41.80400533333333,125.412016,0.9940722981496158,0.9940722981496158,0.6895013292995555,6.205511963695999,0.005179052794781445,0.005179052794781445),
  @("MuSCs","Spn","Siglec1","ECs",2,0.6666666666666666,0.133398,0.400194,0.04213708340045676,0.04213708340045676,2,0.6666666666666666,0.216469,0.6494070000000001,0.005147493274683087,0.005147493274683087,0.028876531662,0.259888784958,0.0002169003534186115,0.0002169003534186115),
  @("MuSCs","Spn","Siglec1","MuSCs",2,0.6666666666666666,0.133398,0.400194,0.04213708340045676,0.04213708340045676,2,0.6666666666666666,0.03281033333333333,0.09843099999999999,0.0007802085757011102,0.0007802085757011101,0.004376832845999999,0.039391495614,0.00003287571382406926,0.00003287571382406926),
  @("MuSCs","Spn","Siglec1","Resolving-Mac",2,0.6666666666666666,0.133398,0.400194,0.04213708340045676,0.04213708340045676,3,1,41.80400533333333,125.412016,0.9940722981496158,0.9940722981496158,5.576570703455999,50.189136331104,0.04188730733321408,0.04188730733321409),
  @("Resolving-Mac","Spn","Siglec1","ECs",3,1,3.015918333333333,9.047754999999999,0.9526529808590324,0.9526529808590325,2,0.6666666666666666,0.216469,0.6494070000000001,0.005147493274683087,0.005147493274683087,0.6528528256983334,5.875675431285,0.004903774812078664,0.004903774812078665),
  @("Resolving-Mac","Spn","Siglec1","MuSCs",3,1,3.015918333333333,9.047754999999999,0.9526529808590324,0.9526529808590325,2,0.6666666666666666,0.03281033333333333,0.09843099999999999,0.0007802085757011102,0.0007802085757011101,0.09895328582277776,0.8905795724049997,0.0007432680253334426,0.0007432680253334426),
  @("Resolving-Mac","Spn","Siglec1","Resolving-Mac",3,1,3.015918333333333,9.047754999999999,0.9526529808590324,0.9526529808590325,3,1,41.80400533333333,125.412016,0.9940722981496158,0.9940722981496158,126.0774660915644,1134.69719482408,0.9470059380216204,0.9470059380216205)
)

for ($r = 0; $r -lt $rows.Length; $r++) {
    $rowVals = $rows[$r]
    for ($c = 0; $c -lt $rowVals.Length; $c++) {
        $ws.Cells.Item($r + 2, $c + 1).Value = $rowVals[$c]
    }
}
